$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "mu" column (I) entirely; everything to its right shifts left.
$ws.Range("I1").EntireColumn.Delete()

# Update each data cell to the new computed values (new column layout A-M).
$ws.Range("C2").Value = 70.95848971701832
$ws.Range("D2").Value = 79561592
$ws.Range("E2").Value = 1126927007.043005
$ws.Range("F2").Value = 0.004220205241383404
$ws.Range("G2").Value = 25.04851854875356
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = $True
$ws.Range("J2").Value = 76923.07692307692
$ws.Range("K2").Value = 109736.707696296
$ws.Range("L2").Value = 18.08138809386727
$ws.Range("M2").Value = 43.12990664262082
$ws.Range("C3").Value = 175.9556230791904
$ws.Range("D3").Value = 79873912
$ws.Range("E3").Value = 232371787.401807
$ws.Range("F3").Value = 0.003491682868916071
$ws.Range("G3").Value = 302.4089996224546
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = $False
$ws.Range("J3").Value = 76923.07692307692
$ws.Range("K3").Value = 97002.55162876946
$ws.Range("L3").Value = 20.45504954955761
$ws.Range("M3").Value = 322.8640491720122
$ws.Range("C4").Value = 118.5748232692446
$ws.Range("D4").Value = 79817080
$ws.Range("E4").Value = 1099128883.62274
$ws.Range("F4").Value = 0.003200154175598993
$ws.Range("G4").Value = 43.0536231732567
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = $True
$ws.Range("J4").Value = 76923.07692307692
$ws.Range("K4").Value = 91467.1420295153
$ws.Range("L4").Value = 21.69294848372683
$ws.Range("M4").Value = 64.74657165698353
$ws.Range("C5").Value = 80.7821316518251
$ws.Range("D5").Value = 79235960
$ws.Range("E5").Value = 833180992.2729816
$ws.Range("F5").Value = 0.0029241020420135
$ws.Range("G5").Value = 38.41212060549256
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = $False
$ws.Range("J5").Value = 76923.07692307692
$ws.Range("K5").Value = 85958.17114104528
$ws.Range("L5").Value = 23.08322726811183
$ws.Range("M5").Value = 61.49534787360439
$ws.Range("C6").Value = 80.74375796642178
$ws.Range("D6").Value = 79761784
$ws.Range("E6").Value = 895042835.8966057
$ws.Range("F6").Value = 0.00430839303353433
$ws.Range("G6").Value = 35.97741875568727
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = $True
$ws.Range("J6").Value = 76923.07692307692
$ws.Range("K6").Value = 111183.5560113591
$ws.Range("L6").Value = 17.84609227462814
$ws.Range("M6").Value = 53.8235110303154
$ws.Range("C7").Value = 111.7396685830363
$ws.Range("D7").Value = 80025976
$ws.Range("E7").Value = 1405892873.14731
$ws.Range("F7").Value = 0.002292505142216274
$ws.Range("G7").Value = 31.80212449706709
$ws.Range("H7").Value = 20
$ws.Range("I7").Value = $True
$ws.Range("J7").Value = 76923.07692307692
$ws.Range("K7").Value = 72219.37576738487
$ws.Range("L7").Value = 27.47451053012403
$ws.Range("M7").Value = 59.27663502719112
$ws.Range("C8").Value = 123.9098299632842
$ws.Range("D8").Value = 79814520
$ws.Range("E8").Value = 1235208275.756846
$ws.Range("F8").Value = 0.0009087328629546882
$ws.Range("G8").Value = 40.03293936701237
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = $False
$ws.Range("J8").Value = 76923.07692307692
$ws.Range("K8").Value = 34408.26517493059
$ws.Range("L8").Value = 57.66614474494507
$ws.Range("M8").Value = 97.69908411195745
$ws.Range("C9").Value = 110.2705246966666
$ws.Range("D9").Value = 79492984
$ws.Range("E9").Value = 967188858.4782858
$ws.Range("F9").Value = 0.002054846491575216
$ws.Range("G9").Value = 45.31551919019817
$ws.Range("H9").Value = 15
$ws.Range("I9").Value = $True
$ws.Range("J9").Value = 76923.07692307692
$ws.Range("K9").Value = 66574.94907215168
$ws.Range("L9").Value = 29.80388310698668
$ws.Range("M9").Value = 75.11940229718485
$ws.Range("C10").Value = 109.2942173926688
$ws.Range("D10").Value = 80046456
$ws.Range("E10").Value = 1828838693.402923
$ws.Range("F10").Value = 0.002065237360215993
$ws.Range("G10").Value = 23.91849755567599
$ws.Range("H10").Value = 20
$ws.Range("I10").Value = $True
$ws.Range("J10").Value = 76923.07692307692
$ws.Range("K10").Value = 66827.82902013171
$ws.Range("L10").Value = 29.69110367781463
$ws.Range("M10").Value = 53.60960123349062
$ws.Range("C11").Value = 153.218164522722
$ws.Range("D11").Value = 79497592
$ws.Range("E11").Value = 788234814.5322708
$ws.Range("F11").Value = 0.002172038545255918
$ws.Range("G11").Value = 77.26425492538016
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = $False
$ws.Range("J11").Value = 76923.07692307692
$ws.Range("K11").Value = 69394.16014641368
$ws.Range("L11").Value = 28.59306886650957
$ws.Range("M11").Value = 105.8573237918897
$ws.Range("C12").Value = 52.69345785456528
$ws.Range("D12").Value = 79762808
$ws.Range("E12").Value = 1304380947.56723
$ws.Range("F12").Value = 0.003557886963972279
$ws.Range("G12").Value = 16.11100717757591
$ws.Range("H12").Value = 20
$ws.Range("I12").Value = $True
$ws.Range("J12").Value = 76923.07692307692
$ws.Range("K12").Value = 98222.0433687905
$ws.Range("L12").Value = 20.20108655803495
$ws.Range("M12").Value = 36.31209373561087
$ws.Range("C13").Value = 185.9418346898805
$ws.Range("D13").Value = 80041848
$ws.Range("E13").Value = 983761988.3588175
$ws.Range("F13").Value = 0.002715096879788233
$ws.Range("G13").Value = 75.64394764793487
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = $False
$ws.Range("J13").Value = 76923.07692307692
$ws.Range("K13").Value = 81597.38144287353
$ws.Range("L13").Value = 24.31685876328196
$ws.Range("M13").Value = 99.96080641121682
$ws.Range("C14").Value = 197.075549015551
$ws.Range("D14").Value = 79997816
$ws.Range("E14").Value = 1645456896.429853
$ws.Range("F14").Value = 0.001254936285397894
$ws.Range("G14").Value = 47.90649193683429
$ws.Range("H14").Value = 9
$ws.Range("I14").Value = $False
$ws.Range("J14").Value = 76923.07692307692
$ws.Range("K14").Value = 45143.10254511117
$ws.Range("L14").Value = 43.95338131705085
$ws.Range("M14").Value = 91.85987325388514
$ws.Range("C15").Value = 72.50100601321395
$ws.Range("D15").Value = 79731064
$ws.Range("E15").Value = 1371947743.31511
$ws.Range("F15").Value = 0.003943618243461979
$ws.Range("G15").Value = 21.06706461186349
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = $True
$ws.Range("J15").Value = 76923.07692307692
$ws.Range("K15").Value = 105072.5329416696
$ws.Range("L15").Value = 18.88402177476308
$ws.Range("M15").Value = 39.95108638662657
$ws.Range("C16").Value = 96.3739219971601
$ws.Range("D16").Value = 79832952
$ws.Range("E16").Value = 789119212.298843
$ws.Range("F16").Value = 0.002360910960038808
$ws.Range("G16").Value = 48.74938139218275
$ws.Range("H16").Value = 15
$ws.Range("I16").Value = $True
$ws.Range("J16").Value = 76923.07692307692
$ws.Range("K16").Value = 73792.200988313
$ws.Range("L16").Value = 26.88891201814472
$ws.Range("M16").Value = 75.63829341032746
$ws.Range("C17").Value = 132.3194100858378
$ws.Range("D17").Value = 79816568
$ws.Range("E17").Value = 1734880216.672365
$ws.Range("F17").Value = 0.003627136141331431
$ws.Range("G17").Value = 30.4380702809947
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = $True
$ws.Range("J17").Value = 76923.07692307692
$ws.Range("K17").Value = 99483.44712401924
$ws.Range("L17").Value = 19.94494619317365
$ws.Range("M17").Value = 50.38301647416836
$ws.Range("C18").Value = 176.1311528522601
$ws.Range("D18").Value = 79846264
$ws.Range("E18").Value = 1118402169.629412
$ws.Range("F18").Value = 0.004926720339789661
$ws.Range("G18").Value = 62.87279706335826
$ws.Range("H18").Value = 13
$ws.Range("I18").Value = $True
$ws.Range("J18").Value = 76923.07692307692
$ws.Range("K18").Value = 120830.5506456392
$ws.Range("L18").Value = 16.42127747823526
$ws.Range("M18").Value = 79.29407454159352
$ws.Range("C19").Value = 102.2661244817619
$ws.Range("D19").Value = 79854968
$ws.Range("E19").Value = 1541271038.75871
$ws.Range("F19").Value = 0.004131978223816866
$ws.Range("G19").Value = 26.49260867365715
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = $True
$ws.Range("J19").Value = 76923.07692307692
$ws.Range("K19").Value = 108270.0903728412
$ws.Range("L19").Value = 18.32631702039958
$ws.Range("M19").Value = 44.81892569405673
$ws.Range("C20").Value = 165.4740402215803
$ws.Range("D20").Value = 80000888
$ws.Range("E20").Value = 1386432394.135247
$ws.Range("F20").Value = 0.002262925582568383
$ws.Range("G20").Value = 47.74149181262842
$ws.Range("H20").Value = 15
$ws.Range("I20").Value = $True
$ws.Range("J20").Value = 76923.07692307692
$ws.Range("K20").Value = 71532.30099033478
$ws.Range("L20").Value = 27.73840590236427
$ws.Range("M20").Value = 75.47989771499269
$ws.Range("C21").Value = 71.22950717745928
$ws.Range("D21").Value = 80230776
$ws.Range("E21").Value = 427332063.0859565
$ws.Range("F21").Value = 0.003963270018734173
$ws.Range("G21").Value = 66.86601742069158
$ws.Range("H21").Value = 9
$ws.Range("I21").Value = $False
$ws.Range("J21").Value = 76923.07692307692
$ws.Range("K21").Value = 105410.4746201164
$ws.Range("L21").Value = 18.82348037185803
$ws.Range("M21").Value = 85.6894977925496
